# Update positionChangeDate (column S) timestamps to reflect the 2024-02-26 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12 get the .504Z timestamp, rows 13-22 get the .505Z timestamp,
# matching the values captured in the updated source export.
for ($row = 2; $row -le 12; $row++) {
    $ws.Range("S$row").Value = "2024-02-26T13:03:26.504Z"
}

for ($row = 13; $row -le 22; $row++) {
    $ws.Range("S$row").Value = "2024-02-26T13:03:26.505Z"
}
